$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) stored its value as a dash-joined string like
# "6-1-2012-13" (month-day-year-of-season-end). The NBA stats export was
# off by one day, so normalize every row to the correct ISO date string
# "2013-06-01" while keeping the cell as plain text (not an Excel date
# serial number).
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # column BF = 58
    $cell.NumberFormat = "@"
    $cell.Value = "2013-06-01"
    $cell.ClearFormats()
}
